$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '42.698.63'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -0.65%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.529.13'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -1.93%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '308.96'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -2.07%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '100.67'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +1.54%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.568'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -1.44%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.07%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.523'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -2.77%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '35.75'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -1.39%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -0.96%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.32'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -3.29%  '
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +0.14%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.916.53'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -1.84%  '
$ws.Range('B15').Value = 'Chainlink'
$ws.Range('C15').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.42'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -1.27%  '
$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.626.05'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +2.39%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.809'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -4.05%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '42.691.82'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -0.85%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.73'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -1.97%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0951'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -2.05%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.26'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -3.23%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '69.42'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -0.05%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '243.63'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -2.82%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.88'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -3.27%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.03'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -3.42%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.01%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '25.50'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -5.92%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.33'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -2.99%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '10.12'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -1.72%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '38.72'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -4.46%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '157.84'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +0.16%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.75'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -1.68%  '
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +11.79%  '
$ws.Range('B34').Value = 'Hedera'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0785'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -2.50%  '
$ws.Range('B35').Value = 'WEMIXToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.64'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -1.58%  '
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -7.92%  '
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -6.78%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '17.77'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -5.39%  '
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -1.23%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -0.75%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.20'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +4.08%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '21.96'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -8.82%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0301'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -1.43%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.28'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +0.39%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.008.37'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +0.10%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.90'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -0.19%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.769.34'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -1.93%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.190'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -3.72%  '
$ws.Range('B50').Value = 'BitcoinSV'
$ws.Range('C50').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '79.25'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -3.50%  '
$ws.Range('B51').Value = 'ordi'
$ws.Range('C51').Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '72.25'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -3.84%  '
